$wb = $excel.ActiveWorkbook

# --- Sheet 1: Cases by Age Group ---
$ws1 = $wb.Worksheets.Item("Cases by Age Group")
$ws1.Range("B2").Value = 262
$ws1.Range("B3").Value = 1251
$ws1.Range("B4").Value = 3376
$ws1.Range("B5").Value = 14683
$ws1.Range("B6").Value = 16200
$ws1.Range("B7").Value = 14164
$ws1.Range("B8").Value = 11971
$ws1.Range("B9").Value = 4320
$ws1.Range("B10").Value = 2888
$ws1.Range("B11").Value = 1714
$ws1.Range("B12").Value = 1103
$ws1.Range("B13").Value = 1716
$ws1.Range("B20").Select()

# --- Sheet 2: Cases by Gender ---
$ws2 = $wb.Worksheets.Item("Cases by Gender")
$ws2.Range("B2").Value = 24724
$ws2.Range("B3").Value = 48028
$ws2.Range("B4").Value = 910

# --- Sheet 3: Cases by RaceEthnicity ---
$ws3 = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws3.Range("B2").Value = 928
$ws3.Range("B3").Value = 12369
$ws3.Range("B4").Value = 27435
$ws3.Range("B5").Value = 400
$ws3.Range("B6").Value = 24275
$ws3.Range("B7").Value = 8255
$ws3.Range("B19").Select()

# --- Sheet 4: Fatalities by Age Group ---
$ws4 = $wb.Worksheets.Item("Fatalities by Age Group")
$ws4.Range("B4").Value = 26
$ws4.Range("B5").Value = 192
$ws4.Range("B6").Value = 632
$ws4.Range("B7").Value = 1867
$ws4.Range("B8").Value = 4320
$ws4.Range("B9").Value = 3659
$ws4.Range("B10").Value = 4689
$ws4.Range("B11").Value = 5323
$ws4.Range("B12").Value = 5351
$ws4.Range("B13").Value = 14022
$ws4.Range("C20").Select()

# --- Sheet 5: Fatalities by Gender ---
$ws5 = $wb.Worksheets.Item("Fatalities by Gender")
$ws5.Range("B2").Value = 16854
$ws5.Range("B3").Value = 23240
$ws5.Range("D14").Select()

# --- Sheet 6: Fatalities by Race-Ethnicity ---
$ws6 = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws6.Range("B2").Value = 800
$ws6.Range("B3").Value = 3816
$ws6.Range("B4").Value = 18648
$ws6.Range("B5").Value = 211
$ws6.Range("B6").Value = 16598
$ws6.Range("B7").Value = 22
$ws6.Range("D17").Select()

# Re-select sheet 1 as active tab (tabSelected)
$ws1.Activate()
$ws1.Range("B20").Select()
